# Insert two new weekly-price rows for "Cilantro" (Vega Central Mapocho de Santiago)
# immediately above the existing row 945. Excel's native row Insert shifts all
# subsequent rows (945-1034) down by two (-> 947-1036) and extends the sheet
# dimension automatically, exactly like inserting rows in the desktop app.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting two blank rows at 945.
$ws.Rows("945:946").Insert()

# ---- New row 945: "$/caja 36 atados" record for the new date ----
$ws.Range("A945").Value = 9
$ws.Range("B945").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C945").Value = "Metropolitana"
$ws.Range("D945").Value = 45132
$ws.Range("E945").Value = 13
$ws.Range("F945").Value = 100112040
$ws.Range("G945").Value = "Cilantro"
$ws.Range("H945").Value = "Sin especificar"
$ws.Range("I945").Value = "Primera"
$ws.Range("J945").Value = 70
$ws.Range("K945").Value = 7000
$ws.Range("L945").Value = 7000
$ws.Range("M945").Value = 7000
$ws.Range("N945").Value = "$/caja 36 atados"
$ws.Range("O945").Value = "Región Metropolitana"
$ws.Range("P945").Value = 194
$ws.Range("Q945").Value = 36
$ws.Range("R945").Value = "Hortaliza"

# ---- New row 946: "$/docena de atados" record for the same new date ----
$ws.Range("A946").Value = 9
$ws.Range("B946").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C946").Value = "Metropolitana"
$ws.Range("D946").Value = 45132
$ws.Range("E946").Value = 13
$ws.Range("F946").Value = 100112040
$ws.Range("G946").Value = "Cilantro"
$ws.Range("H946").Value = "Sin especificar"
$ws.Range("I946").Value = "Primera"
$ws.Range("J946").Value = 160
$ws.Range("K946").Value = 11000
$ws.Range("L946").Value = 12000
$ws.Range("M946").Value = 11500
$ws.Range("N946").Value = "$/docena de atados"
$ws.Range("O946").Value = "Región Metropolitana"
$ws.Range("P946").Value = 3833
$ws.Range("Q946").Value = 3
$ws.Range("R946").Value = "Hortaliza"
